$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A16").Value = 14

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 0.9859038466344734
$ws.Range("D16").Value = 1.017880101365173
$ws.Range("E16").Value = 0.9888596471693042
$ws.Range("F16").Value = 0.9859038466344734
$ws.Range("G16").Value = 1.008540572174369
$ws.Range("H16").Value = 0.981394061939606
$ws.Range("I16").Value = 0.9882060459891621
$ws.Range("J16").Value = 1.017880101365173
$ws.Range("K16").Value = 1.003369874267239
$ws.Range("L16").Value = 0.994636860450856
$ws.Range("M16").Value = 0.995130712545348
